$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3934
$ws1.Range("F4").Value = 2319
$ws1.Range("F5").Value = 459
$ws1.Range("F10").Value = 15
$ws1.Range("F11").Value = 116
$ws1.Range("F14").Value = 2651
$ws1.Range("F15").Value = 182

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3934
$ws4.Range("F4").Value = 2319
$ws4.Range("F5").Value = 459
$ws4.Range("F11").Value = 15
$ws4.Range("F12").Value = 116
$ws4.Range("F17").Value = 2651
$ws4.Range("F18").Value = 182
